# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted as row 151 of the data table
# (the sheet's only worksheet). Inserting the row pushes the existing
# rows 151..235 down to 152..236, which is exactly what the target
# workbook's diff shows (every row below 150 keeps its original values,
# just shifted down by one row number; the sheet dimension grows from
# A1:R235 to A1:R236).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at row 151 - shifts rows 151:235 -> 152:236
# and grows the used range / dimension automatically.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new weekly record.
$ws.Cells.Item(151, 1).Value  = 3
$ws.Cells.Item(151, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(151, 3).Value  = "Coquimbo"
$ws.Cells.Item(151, 4).Value  = 44518
$ws.Cells.Item(151, 5).Value  = 5
$ws.Cells.Item(151, 6).Value  = 100112043
$ws.Cells.Item(151, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(151, 8).Value  = "Sin especificar"
$ws.Cells.Item(151, 9).Value  = "Primera"
$ws.Cells.Item(151, 10).Value = 125
$ws.Cells.Item(151, 11).Value = 7000
$ws.Cells.Item(151, 12).Value = 8000
$ws.Cells.Item(151, 13).Value = 7520
$ws.Cells.Item(151, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(151, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(151, 16).Value = 107
$ws.Cells.Item(151, 17).Value = 70
$ws.Cells.Item(151, 18).Value = "Hortaliza"
